$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 47210.363
$ws.Range("J112").Value = 60956.41
$ws.Range("L112").Value = 182869.23
$ws.Range("N112").Value = -185085.23
$ws.Range("H116").Value = 2353.2354
$ws.Range("I116").Value = 1750.8334
$ws.Range("J116").Value = 2681.818
$ws.Range("K116").Value = 1750.8334
$ws.Range("L116").Value = 2681.818
$ws.Range("M116").Value = 1691.1666
$ws.Range("N116").Value = -9565.817999999999
$ws.Range("H127").Value = 1071
$ws.Range("I127").Value = 199.1
$ws.Range("J127").Value = 2160.875
$ws.Range("K127").Value = 597.3
$ws.Range("L127").Value = 6482.625
$ws.Range("M127").Value = 4362.7
$ws.Range("N127").Value = -16402.625
$ws.Range("H138").Value = 2751.0527
$ws.Range("I138").Value = 1101.5
$ws.Range("J138").Value = 3826.848
$ws.Range("K138").Value = 3304.5
$ws.Range("L138").Value = 11480.544
$ws.Range("M138").Value = 1835.5
$ws.Range("N138").Value = -21760.544
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4832.8
$ws.Range("I32").Value = 3688.375
$ws.Range("J32").Value = 9410.5
$ws.Range("K32").Value = 3688.375
$ws.Range("L32").Value = 9410.5
$ws.Range("M32").Value = -3401.375
$ws.Range("N32").Value = -9984.5
$ws.Range("H45").Value = 1922.2778
$ws.Range("I45").Value = 1750.1428
$ws.Range("J45").Value = 2524.75
$ws.Range("K45").Value = 1750.1428
$ws.Range("L45").Value = 2524.75
$ws.Range("M45").Value = -1373.1428
$ws.Range("N45").Value = -3278.75
$ws.Range("H122").Value = 2631.111
$ws.Range("I122").Value = 2525
$ws.Range("J122").Value = 2716
$ws.Range("K122").Value = 7575
$ws.Range("L122").Value = 8148
$ws.Range("M122").Value = -5125
$ws.Range("N122").Value = -13048
$ws.Range("H132").Value = 23191.555
$ws.Range("I132").Value = 17220.834
$ws.Range("J132").Value = 35133
$ws.Range("K132").Value = 51662.50199999999
$ws.Range("L132").Value = 105399
$ws.Range("M132").Value = -49132.50199999999
$ws.Range("N132").Value = -110459
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 12271.2
$ws.Range("I96").Value = 12271.2
$ws.Range("K96").Value = 12271.2
$ws.Range("M96").Value = -9525.200000000001
$ws.Range("H134").Value = 6665.294
$ws.Range("I134").Value = 6665.294
$ws.Range("K134").Value = 19995.882
$ws.Range("M134").Value = -17460.882
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3147
$ws.Range("J31").Value = 3147
$ws.Range("L31").Value = 3147
$ws.Range("N31").Value = -3737
$ws.Range("H34").Value = 3147
$ws.Range("J34").Value = 3147
$ws.Range("L34").Value = 3147
$ws.Range("N34").Value = -3551
$ws.Range("H99").Value = 2037.6666
$ws.Range("I99").Value = 1804
$ws.Range("J99").Value = 2271.3333
$ws.Range("K99").Value = 1804
$ws.Range("L99").Value = 2271.3333
$ws.Range("M99").Value = -306
$ws.Range("N99").Value = -5267.3333
$ws.Range("H122").Value = 2541.1538
$ws.Range("I122").Value = 3240.2
$ws.Range("J122").Value = 2104.25
$ws.Range("K122").Value = 9720.599999999999
$ws.Range("L122").Value = 6312.75
$ws.Range("M122").Value = -7270.599999999999
$ws.Range("N122").Value = -11212.75
$ws.Range("H126").Value = 2037.6666
$ws.Range("I126").Value = 1804
$ws.Range("J126").Value = 2271.3333
$ws.Range("K126").Value = 5412
$ws.Range("L126").Value = 6813.999899999999
$ws.Range("M126").Value = -2942
$ws.Range("N126").Value = -11753.9999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2435.96
$ws.Range("I5").Value = 519.2
$ws.Range("J5").Value = 2915.15
$ws.Range("K5").Value = 1557.6
$ws.Range("L5").Value = 8745.450000000001
$ws.Range("M5").Value = -1445.6
$ws.Range("N5").Value = -8969.450000000001
$ws.Range("H122").Value = 14257.643
$ws.Range("I122").Value = 304
$ws.Range("J122").Value = 15331
$ws.Range("K122").Value = 2736
$ws.Range("L122").Value = 137979
$ws.Range("M122").Value = -286
$ws.Range("N122").Value = -142879
$ws.Range("H132").Value = 1314.2941
$ws.Range("I132").Value = 905.4286
$ws.Range("J132").Value = 1600.5
$ws.Range("K132").Value = 8148.8574
$ws.Range("L132").Value = 14404.5
$ws.Range("M132").Value = -5618.8574
$ws.Range("N132").Value = -19464.5
$ws.Range("H134").Value = 3969.8096
$ws.Range("I134").Value = 1738.4615
$ws.Range("J134").Value = 7595.75
$ws.Range("K134").Value = 5215.3845
$ws.Range("L134").Value = 22787.25
$ws.Range("M134").Value = -145.3845000000001
$ws.Range("N134").Value = -32927.25
$ws.Range("H135").Value = 2435.96
$ws.Range("I135").Value = 519.2
$ws.Range("J135").Value = 2915.15
$ws.Range("K135").Value = 4672.8
$ws.Range("L135").Value = 26236.35
$ws.Range("M135").Value = -2137.8
$ws.Range("N135").Value = -31306.35
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 25333.111
$ws.Range("I132").Value = 19333.334
$ws.Range("K132").Value = 58000.00199999999
$ws.Range("M132").Value = -55470.00199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1636.1428
$ws.Range("I61").Value = 1490.3889
$ws.Range("J61").Value = 1898.5
$ws.Range("K61").Value = 1490.3889
$ws.Range("L61").Value = 1898.5
$ws.Range("M61").Value = -1288.3889
$ws.Range("N61").Value = -2302.5
$ws.Range("H82").Value = 2729.6
$ws.Range("I82").Value = 1866.3334
$ws.Range("J82").Value = 2881.9412
$ws.Range("K82").Value = 1866.3334
$ws.Range("L82").Value = 2881.9412
$ws.Range("M82").Value = -1505.3334
$ws.Range("N82").Value = -3603.9412
$ws.Range("H85").Value = 2729.6
$ws.Range("I85").Value = 1866.3334
$ws.Range("J85").Value = 2881.9412
$ws.Range("K85").Value = 1866.3334
$ws.Range("L85").Value = 2881.9412
$ws.Range("M85").Value = -618.3334
$ws.Range("N85").Value = -5377.9412
$ws.Range("H100").Value = 2661.6858
$ws.Range("I100").Value = 1584.875
$ws.Range("J100").Value = 2980.7407
$ws.Range("K100").Value = 1584.875
$ws.Range("L100").Value = 2980.7407
$ws.Range("M100").Value = -1043.875
$ws.Range("N100").Value = -4062.7407
$ws.Range("H113").Value = 1636.1428
$ws.Range("I113").Value = 1490.3889
$ws.Range("J113").Value = 1898.5
$ws.Range("K113").Value = 1490.3889
$ws.Range("L113").Value = 1898.5
$ws.Range("M113").Value = 679.6111000000001
$ws.Range("N113").Value = -6238.5
$ws.Range("H122").Value = 2738.8
$ws.Range("I122").Value = 2787.25
$ws.Range("J122").Value = 2545
$ws.Range("K122").Value = 8361.75
$ws.Range("L122").Value = 7635
$ws.Range("M122").Value = -5911.75
$ws.Range("N122").Value = -12535
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1918194.8
$ws.Range("I122").Value = 1059693.8
$ws.Range("J122").Value = 3791288
$ws.Range("K122").Value = 3179081.4
$ws.Range("L122").Value = 11373864
$ws.Range("M122").Value = -3176631.4
$ws.Range("N122").Value = -11378764
$ws.Range("H126").Value = 3307417.8
$ws.Range("I126").Value = 3980134.2
$ws.Range("K126").Value = 11940402.6
$ws.Range("M126").Value = -11937932.6
